# Generate Report for Handoff
# Updates the "latest handoff" timestamps for the 51ca10d6-... file row
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-06 13:52:57"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for row 4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-06 13:52:46"

# de-de sheet: "Latest Handoff Datetime" column (H) for row 4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-06 13:52:57"
